# feat: table foot sum D and C
# Apply the highlight fill (same style already used by rows 2-6, cellXf index 2
# -> fill color D8E4BC) to rows 7,8,9,11,13,16,17 across columns A:K, and flip
# the "conciliada" (K) flag from 0 to 1 for those same rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Interior.Color read from an already-styled row (row 2, which carries the
# fill used throughout this table) expressed as an OLE BGR value.
$fillColor = $ws.Range("A2").Interior.Color

$rows = @(7, 8, 9, 11, 13, 16, 17)

foreach ($r in $rows) {
    $rowRange = $ws.Range("A" + $r + ":K" + $r)
    $rowRange.Interior.Color = $fillColor

    $kCell = $ws.Cells.Item($r, 11)
    $kCell.Value = 1
}
